# Applies the "Suppletion SourceDatatype and correction EventDate" edit:
#  - Fills in column G (SourceDatatypeEng) for rows 2-78 with datatype labels
#  - Adjusts a handful of view/workbook cosmetic properties
#
# Column G values are copied in (value + matching cell format) using donor
# cells that already carry the exact target format, via Copy/PasteSpecial,
# so that the resulting style indexes line up with the ones Excel itself
# would reuse (no new, redundant style entries are introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Donor cells already present on the sheet that carry each of the cell
# formats we need to reproduce for column G (picked before any edits are
# made, so they are unaffected by the changes below).
# ---------------------------------------------------------------------------
$styleDonor = @{
    1  = "R12"   # vertical=top
    2  = "A5"    # border bottom (thin)
    3  = "R5"    # border bottom (thin) + vertical=top
    6  = "J2"    # horizontal=left
    7  = "R2"    # vertical=top + wrap text
    9  = "K38"   # horizontal=left + vertical=top
    12 = "R78"   # border bottom (medium) + vertical=top
}

function Set-GCell($row, $value, $styleId) {
    $target = $ws.Range("G$row")
    if ($styleId -eq $null) {
        $target.ClearFormats()
        $target.Value = $value
    } else {
        $donor = $ws.Range($styleDonor[$styleId])
        $donor.Copy()
        $target.PasteSpecial(-4122)  # xlPasteFormats
        $target.Value = $value
    }
}

# ---------------------------------------------------------------------------
# Column G (SourceDatatypeEng) suppletion, rows 2-78
# ---------------------------------------------------------------------------
Set-GCell 2  "Geen"              $null
Set-GCell 3  "Geen"              $null
Set-GCell 4  "Geen"              $null
Set-GCell 5  "[Uint8]"           3
Set-GCell 6  "[BroID]"           $null
Set-GCell 7  "[KvKNumber]"       6
Set-GCell 8  "[String]"          6
Set-GCell 9  "[KvKNumber]"       6
Set-GCell 10 "[IDCategorical]"   2
Set-GCell 11 "[DateTime]"        $null
Set-GCell 12 "[Integer]"         1
Set-GCell 13 "[Categorical]"     $null
Set-GCell 14 "[DateTime]"        $null
Set-GCell 15 "[DateTime]"        $null
Set-GCell 16 "[Categorical]"     $null
Set-GCell 17 "[DateTime]"        $null
Set-GCell 18 "[Categorical]"     $null
Set-GCell 19 "[DateTime]"        $null
Set-GCell 20 "[Categorical]"     $null
Set-GCell 21 "[DateTime]"        $null
Set-GCell 22 "[Categorical]"     $null
Set-GCell 23 "[DateTime]"        3
Set-GCell 24 "[JaNee]"           1
Set-GCell 25 "[PutCode]"         1
Set-GCell 26 "[NITGCode]"        1
Set-GCell 27 "[KvKNumber]"       6
Set-GCell 28 "[KvKNumber]"       6
Set-GCell 29 "[Categorical]"     1
Set-GCell 30 "[Categorical]"     1
Set-GCell 31 "[Categorical]"     1
Set-GCell 32 "[JaNee]"           1
Set-GCell 33 "[Integer]"         1
Set-GCell 34 "[JaNeeOnbekend]"   1
Set-GCell 35 "[Categorical]"     1
Set-GCell 36 "[Categorical]"     3
Set-GCell 37 "[DateTime]"        7
Set-GCell 38 "[DateTime]"        7
Set-GCell 39 "[Categorical]"     1
Set-GCell 40 "[DateTime]"        7
Set-GCell 41 "[CoördinatePair]"  9
Set-GCell 42 "[Categorical]"     1
Set-GCell 43 "[Categorical]"     1
Set-GCell 44 "[Categorical]"     1
Set-GCell 45 "[m]"               1
Set-GCell 46 "[Categorical]"     1
Set-GCell 47 "[m+NAP]"           $null
Set-GCell 48 "[Categorical]"     3
Set-GCell 49 "[Uint8]"           1
Set-GCell 50 "[Categorical]"     1
Set-GCell 51 "[JaNeeOnbekend]"   1
Set-GCell 52 "[JaNeeOnbekend]"   1
Set-GCell 53 "[Integer]"         1
Set-GCell 54 "[JaNee]"           1
Set-GCell 55 "[mm]"              1
Set-GCell 56 "[JaNeeOnbekend]"   1
Set-GCell 57 "[Categorical]"     1
Set-GCell 58 "[m+NAP]"           $null
Set-GCell 59 "[Categorical]"     1
Set-GCell 60 "[JaNeeOnbekend]"   1
Set-GCell 61 "[Integer]"         1
Set-GCell 62 "[JaNeeOnbekend]"   1
Set-GCell 63 "[Integer]"         1
Set-GCell 64 "[Categorical]"     1
Set-GCell 65 "[m]"               1
Set-GCell 66 "[Categorical]"     3
Set-GCell 67 "[Categorical]"     1
Set-GCell 68 "[Categorical]"     1
Set-GCell 69 "[Categorical]"     1
Set-GCell 70 "[m]"               1
Set-GCell 71 "[Categorical]"     1
Set-GCell 72 "[m+NAP]"           $null
Set-GCell 73 "[m+NAP]"           $null
Set-GCell 74 "[m]"               1
Set-GCell 75 "[m]"               1
Set-GCell 76 "[m]"               1
Set-GCell 77 "[m]"               1
Set-GCell 78 "[Categorical]"     12

$excel.CutCopyMode = 0
